# Update cryptocurrency price/volume data scraped on Fri Apr 21 04:26:15 UTC 2023.
# For D-column cells whose new value looks like a plain decimal number (e.g. "1.009"),
# temporarily force a Text number format while assigning the value so Excel keeps the
# original text representation (matching the source inlineStr cells) instead of silently
# converting it to a floating point number, then restore the default "Normal" style so the
# cell's style index is unchanged from the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.553.93'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '1.964.84'
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4834'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4083'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.94'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("E10").Value = '  -6.38%  '
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.55'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("D13").Value = '2.018.56'
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.632'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.00%  '
$ws.Range("E15").Value = '  -3.47%  '
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06629'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("E20").Value = '  -3.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.010'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.886'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("D23").Value = '28.613.45'
$ws.Range("E23").Value = '  -1.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.298'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("D26").Value = '2.271.66'
$ws.Range("E26").Value = '  +3.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.931'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("E30").Value = '  -3.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '125.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9966'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09735'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.474'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.683'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.696'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.180'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02344'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06266'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.260'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6265'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.30'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.009'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1926'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.351'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.88%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5988'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.57%  '
$ws.Range("E48").Value = '  -4.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.414'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06840'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.63'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.84%  '
